$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'279.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.52%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'27.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.89%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.832"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.91%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06384"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.24%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.040"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.93%"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'5.02%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8955"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.67%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1546"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.61%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.06567"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'28.25%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07540"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.36%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.02945"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.68%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09004"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.06%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001566"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.10%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006450"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.86%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006065"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'2.99%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.487"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.74%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'0.22%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.230"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.86%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D21").Value = "'0.1351"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.06%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.910"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.10%"
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "'0.1504"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'8.97%"
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "CoinExToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D24").Value = "'0.04399"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.46%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001176"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.27%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004282"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'10.36%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E28").Value = "'-1.66%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'0.0001654"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'-14.56%"
$ws.Range("E29").Style = "Normal"
$ws.Range("D40").Value = "'0.04067"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.94%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006641"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-2.57%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1411"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'19.62%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002090"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'3.47%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'-1.79%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005556"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'7.11%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.628"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'9.52%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.01850"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-8.62%"
$ws.Range("E47").Style = "Normal"
